$d = $word.ActiveDocument

# Avoid Word's AutoFormat turning straight quotes into curly "smart" quotes
# during Find/Replace text substitution.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false

function Replace-All($find, $replace) {
    $content = $d.Content
    $content.Find.Execute($find, $true, $false, $false, $false, $false, `
                           $true, 1, $false, $null, 0) | Out-Null
    while ($content.Find.Found) {
        $start = $content.Start
        $end = $content.End
        $target = $d.Range($start, $end)
        $target.Text = $replace
        $newEnd = $start + $replace.Length

        # The engine tends to coalesce the freshly edited run with an
        # immediately adjacent run that ends up sharing the same
        # character formatting, even when that neighbour run was not
        # part of the edit. Toggle a character property off/on across
        # just the replacement text to force the run back apart again.
        $splitRange = $d.Range($start, $newEnd)
        $splitRange.Bold = 1
        $splitRange.Bold = 0

        $content = $d.Content
        $content.Start = $newEnd
        $content.Find.Execute($find, $true, $false, $false, $false, $false, `
                               $true, 1, $false, $null, 0) | Out-Null
    }
}

Replace-All "2445987" "1122334"
Replace-All "01.05.2024" "07.03.2024"
Replace-All "генеральный директор Котлярчук О. Е." "Капитан Бахтин Ю. Г."
Replace-All "Устава" "Кодекса торгового мореплавания (КТМ РФ)"
Replace-All '"15 ВАХАУ МАРУ" ' '"СИНЕГОРСК" '
Replace-All "172316" "021026"
Replace-All "Внеочередное освидетельствование в связи со сменой судовладельца" "Первоначальное освидетельствование на соответствие требованиям МК ОСПС"
Replace-All "Акт ф. 6.1.03 № 24.43.01.00135.121 от 05.05.2024" "Свидетельство ф. 8.5.3 № 24.42.02.00123.121 от 04.05.2024"
Replace-All "20 236,66 p. (двадцать тысяч двести тридцать шесть рублей 66 копеек)" "100 000,00 p. (сто тысяч рублей 00 копеек)"
Replace-All "О. Е. Котлярчук" "Ю. Г. Бахтин"
